# chartink_screener.xlsx - "break out stock.yaml completed"
#
# For each of the three scraper-result sheets ("10per change", "3 V 0.3",
# "DND 3 V 0.3") a new scrape run appended a duplicate of the existing
# rows (same data) stamped with a newer "Date Time", and the bsecode
# (column D) on the *original* rows got normalised from text to a real
# number. The "3 V 0.4" sheet is untouched by this run.

$wb = $excel.ActiveWorkbook
$newTimestamp = "02/06/2024 06:03:37"

function Append-ScrapeRun {
    param(
        [string]$SheetName,
        [int]$FirstDataRow,
        [int]$LastDataRow
    )

    $ws = $wb.Worksheets.Item($SheetName)
    $rowCount = $LastDataRow - $FirstDataRow + 1

    $srcRange = $ws.Range("A" + $FirstDataRow + ":H" + $LastDataRow)
    $newFirst = $LastDataRow + 1
    $newLast = $LastDataRow + $rowCount
    $dstRange = $ws.Range("A" + $newFirst + ":H" + $newLast)

    # Duplicate the block of rows verbatim (this also preserves the
    # text/number type of every cell, including column D's text bsecode).
    $srcRange.Copy($dstRange)

    # Stamp the freshly appended rows with the new scrape timestamp.
    $ws.Range("H" + $newFirst + ":H" + $newLast).Value = $newTimestamp

    # Normalise the *original* rows' bsecode (column D) from text to number.
    for ($r = $FirstDataRow; $r -le $LastDataRow; $r++) {
        $cell = $ws.Cells.Item($r, 4)
        $cell.Value = [double]$cell.Text
    }
}

Append-ScrapeRun "10per change" 2 7
Append-ScrapeRun "3 V 0.3" 2 3
Append-ScrapeRun "DND 3 V 0.3" 2 2
